$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CHAM_CONG_HE_THONG")

# Update last_edited_time (column D) for all data rows (2-20) from
# 2024-07-18T15:58:00.000Z to 2024-07-19T08:01:00.000Z
for ($r = 2; $r -le 20; $r++) {
    $ws.Cells.Item($r, 4).Value = "2024-07-19T08:01:00.000Z"
}
